$wb = $excel.ActiveWorkbook

$wsInfo   = $wb.Worksheets.Item("Info")
$ws10be1  = $wb.Worksheets.Item("10Be_Model1")
$ws36cl1  = $wb.Worksheets.Item("36Cl_Model1")
$ws36cl2  = $wb.Worksheets.Item("36Cl_Model2")

# --- 10Be_Model1: update depth/age values ---
$ws10be1.Range("D2").Value = 68.5
$ws10be1.Range("E2").Value = 12.6
$ws10be1.Range("D3").Value = 53.8
$ws10be1.Range("E3").Value = 11.1
$ws10be1.Range("D4").Value = 56.3
$ws10be1.Range("E4").Value = 10.2

# --- 36Cl_Model1: update depth/age values, widen column C ---
$ws36cl1.Range("D2").Value = 13
$ws36cl1.Range("E2").Value = 1.7
$ws36cl1.Range("D4").Value = 13
$ws36cl1.Range("E4").Value = 1.7
$ws36cl1.Columns.Item(3).ColumnWidth = 11.33

# --- Update the active cell / selection on each sheet ---
$ws10be1.Activate() | Out-Null
$ws10be1.Range("C37").Select() | Out-Null

$ws36cl1.Activate() | Out-Null
$ws36cl1.Range("F4").Select() | Out-Null

$ws36cl2.Activate() | Out-Null
$ws36cl2.Range("F10").Select() | Out-Null

# Info stays the tab that is active/selected when the file is saved
$wsInfo.Activate() | Out-Null
$wsInfo.Range("C9").Select() | Out-Null
